# Updated symbol list (price/volume refresh + row reshuffle for HotbitToken/LEO/...)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force these cells to remain text (they hold numeric-looking strings like
# "321.19" or "7.92%") so Excel doesn't silently convert them to numbers.
foreach ($addr in @('D2','E2','D3','E3','D4','E4','D5','E5','D6','E6','D7','E7','D8','E8','D9','E9','D10','E10','D11','E11','D12','E12','D13','E13','D14','E14','D15','E15','B16','C16','D16','E16','B17','C17','D17','E17','B18','C18','D18','E18','B19','C19','D19','E19','B20','C20','D20','E20','B21','C21','D21','E21','B22','C22','D22','E22','B23','C23','D23','E23','B24','C24','D24','E24','D25','D26','E26','D38','E38','D39','E39','D40','E40','D41','E41','D42','E42','D43','E43','D44','E44','D45','E45','D46','E46','D47','E47','E48','D49','E49','D50','E50','D51','E51')) {
    $ws.Range($addr).NumberFormat = '@'
}

$ws.Range('D2').Value = '321.19'
$ws.Range('E2').Value = '7.92%'
$ws.Range('D3').Value = '50.13'
$ws.Range('E3').Value = '19.62%'
$ws.Range('D4').Value = '5.349'
$ws.Range('E4').Value = '6.94%'
$ws.Range('D5').Value = '0.08146'
$ws.Range('E5').Value = '8.33%'
$ws.Range('D6').Value = '4.604'
$ws.Range('E6').Value = '5.05%'
$ws.Range('D7').Value = '1.669'
$ws.Range('E7').Value = '5.37%'
$ws.Range('D8').Value = '1.167'
$ws.Range('E8').Value = '25.91%'
$ws.Range('D9').Value = '0.1332'
$ws.Range('E9').Value = '11.39%'
$ws.Range('D10').Value = '0.1950'
$ws.Range('E10').Value = '6.29%'
$ws.Range('D11').Value = '0.09545'
$ws.Range('E11').Value = '7.10%'
$ws.Range('D12').Value = '0.04574'
$ws.Range('E12').Value = '12.20%'
$ws.Range('D13').Value = '0.1047'
$ws.Range('E13').Value = '-0.05%'
$ws.Range('D14').Value = '0.001331'
$ws.Range('E14').Value = '3.62%'
$ws.Range('D15').Value = '0.005863'
$ws.Range('E15').Value = '-1.74%'
$ws.Range('B16').Value = 'HotbitToken'
$ws.Range('C16').Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
$ws.Range('D16').Value = '0.004300'
$ws.Range('E16').Value = '10.58%'
$ws.Range('B17').Value = 'LEO'
$ws.Range('C17').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D17').Value = '3.380'
$ws.Range('E17').Value = '0.50%'
$ws.Range('B18').Value = 'BTSEToken'
$ws.Range('C18').Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range('D18').Value = '2.432'
$ws.Range('E18').Value = '1.30%'
$ws.Range('B19').Value = 'BitpandaEcosystemToken'
$ws.Range('C19').Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range('D19').Value = '0.3394'
$ws.Range('E19').Value = '2.43%'
$ws.Range('B20').Value = 'MCDex'
$ws.Range('C20').Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range('D20').Value = '8.159'
$ws.Range('E20').Value = '0.70%'
$ws.Range('B21').Value = 'ProBitToken'
$ws.Range('C21').Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$ws.Range('D21').Value = '0.1409'
$ws.Range('E21').Value = '1.35%'
$ws.Range('B22').Value = 'ZBToken'
$ws.Range('C22').Value = 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'
$ws.Range('D22').Value = '0.3055'
$ws.Range('E22').Value = '-7.53%'
$ws.Range('B23').Value = 'CoinExToken'
$ws.Range('C23').Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range('D23').Value = '0.04311'
$ws.Range('E23').Value = '5.12%'
$ws.Range('B24').Value = 'BitKan'
$ws.Range('C24').Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
$ws.Range('D24').Value = '0.001304'
$ws.Range('E24').Value = '3.06%'
$ws.Range('D25').Value = '0.0001349'
$ws.Range('D26').Value = '0.0003720'
$ws.Range('E26').Value = '-0.10%'
$ws.Range('D38').Value = '0.02753'
$ws.Range('E38').Value = '13.96%'
$ws.Range('D39').Value = '0.05524'
$ws.Range('E39').Value = '5.94%'
$ws.Range('D40').Value = '0.006198'
$ws.Range('E40').Value = '-1.68%'
$ws.Range('D41').Value = '0.007775'
$ws.Range('E41').Value = '-0.62%'
$ws.Range('D42').Value = '0.1445'
$ws.Range('E42').Value = '9.01%'
$ws.Range('D43').Value = '0.007690'
$ws.Range('E43').Value = '4.04%'
$ws.Range('D44').Value = '0.008850'
$ws.Range('E44').Value = '15.13%'
$ws.Range('D45').Value = '0.3485'
$ws.Range('E45').Value = '17.79%'
$ws.Range('D46').Value = '0.00006768'
$ws.Range('E46').Value = '5.72%'
$ws.Range('D47').Value = '0.00000000750'
$ws.Range('E47').Value = '-0.10%'
$ws.Range('E48').Value = '96.52%'
$ws.Range('D49').Value = '0.003998'
$ws.Range('E49').Value = '-4.85%'
$ws.Range('D50').Value = '0.00002099'
$ws.Range('E50').Value = '-0.10%'
$ws.Range('D51').Value = '0.0001999'
$ws.Range('E51').Value = '-0.10%'
